# Apply the cryptos list update (cell-level edits from the XML diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E are plain (non-formula) text cells holding formatted
# price/volume strings. Some new D values look like plain numbers to Excel
# (e.g. "167.29"), so the affected D cells are pre-formatted as Text to keep
# them stored verbatim instead of being coerced into numeric values.
$textFormatCells = @("D5", "D6", "D8", "D10", "D13", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "67.174.54"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").Value = "2.614.22"
$ws.Range("E3").Value = "  -2.32%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "594.61"
$ws.Range("E5").Value = "  -0.57%  "

# Row 6
$ws.Range("D6").Value = "167.29"
$ws.Range("E6").Value = "  +1.00%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  -2.32%  "

# Row 9
$ws.Range("D9").Value = "2.613.19"
$ws.Range("E9").Value = "  -2.33%  "

# Row 10
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -2.15%  "

# Row 11
$ws.Range("E11").Value = "  +1.09%  "

# Row 12
$ws.Range("E12").Value = "  +1.24%  "

# Row 13
$ws.Range("D13").Value = "5.22"
$ws.Range("E13").Value = "  +0.04%  "

# Row 14
$ws.Range("D14").Value = "27.53"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15
$ws.Range("D15").Value = "3.091.68"
$ws.Range("E15").Value = "  -2.32%  "

# Row 16
$ws.Range("D16").Value = "0.0000182"
$ws.Range("E16").Value = "  -0.90%  "

# Row 17
$ws.Range("D17").Value = "67.188.75"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18
$ws.Range("D18").Value = "2.626.22"
$ws.Range("E18").Value = "  -1.28%  "

# Row 19
$ws.Range("D19").Value = "12.02"
$ws.Range("E19").Value = "  +2.59%  "

# Row 20
$ws.Range("D20").Value = "7.94"
$ws.Range("E20").Value = "  +4.39%  "

# Row 21
$ws.Range("D21").Value = "356.02"
$ws.Range("E21").Value = "  -1.99%  "

# Row 22
$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  -1.43%  "

# Row 23
$ws.Range("D23").Value = "4.65"
$ws.Range("E23").Value = "  -3.43%  "

# Row 25
$ws.Range("D25").Value = "1.92"
$ws.Range("E25").Value = "  -5.11%  "

# Row 26
$ws.Range("D26").Value = "10.27"
$ws.Range("E26").Value = "  +0.85%  "

# Row 27
$ws.Range("D27").Value = "69.66"
$ws.Range("E27").Value = "  -2.19%  "

# Row 28
$ws.Range("E28").Value = "  -2.00%  "

# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0999"
$ws.Range("E30").Value = "  -2.17%  "

# Row 31
$ws.Range("D31").Value = "541.82"
$ws.Range("E31").Value = "  -2.22%  "

# Row 32
$ws.Range("D32").Value = "7.88"
$ws.Range("E32").Value = "  -1.48%  "

# Row 33
$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  -3.09%  "

# Row 34
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  -1.74%  "

# Row 35
$ws.Range("E35").Value = "  +4.45%  "

# Row 36
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("D37").Value = "1.50"
$ws.Range("E37").Value = "  -4.30%  "

# Row 38
$ws.Range("D38").Value = "157.01"
$ws.Range("E38").Value = "  +1.63%  "

# Row 39
$ws.Range("D39").Value = "18.93"
$ws.Range("E39").Value = "  -2.95%  "

# Row 40
$ws.Range("D40").Value = "0.366"
$ws.Range("E40").Value = "  -2.06%  "

# Row 41
$ws.Range("D41").Value = "18.16"
$ws.Range("E41").Value = "  +1.34%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.80"
$ws.Range("E42").Value = "  -1.33%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "5.19"
$ws.Range("E43").Value = "  -2.04%  "

# Row 44
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  -4.62%  "

# Row 46
$ws.Range("D46").Value = "0.0₆0296"
$ws.Range("E46").Value = "  -0.70%  "

# Row 47
$ws.Range("D47").Value = "151.65"
$ws.Range("E47").Value = "  -0.88%  "

# Row 48
$ws.Range("D48").Value = "0.577"
$ws.Range("E48").Value = "  -2.79%  "

# Row 49
$ws.Range("D49").Value = "3.76"
$ws.Range("E49").Value = "  -1.80%  "

# Row 50
$ws.Range("D50").Value = "1.69"

# Row 51
$ws.Range("D51").Value = "0.0769"
$ws.Range("E51").Value = "  -1.09%  "
